# Add a new "Compromissadas" asset row for each of the three funds that
# already had Stocks / LFT / CLCD16 / PETR4 rows, keeping the funds grouped
# together (new row inserted right after each fund's existing block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new row in the "TREND DI..." block (was rows 2-5) and
# in the "EQUITAS SHELTER" block (shifts down to rows 7-10 after the first
# insert). The last new row ("DRYS SHELTER PREV" block) lands right after
# the current last row of data, so no insert is required for it - it is
# simply written into the next empty row.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(11).Insert()

$ws.Range("A6").Value = "TREND DI SIMPLES FUNDOS DE INVESTIMENTO EM RENDA FIXA"
$ws.Range("B6").Value = "Compromissadas"

$ws.Range("A11").Value = "EQUITAS SHELTER"
$ws.Range("B11").Value = "Compromissadas"

$ws.Range("A16").Value = "DRYS SHELTER PREV"
$ws.Range("B16").Value = "Compromissadas"

# Reflect the cell that was left selected after the edit.
$ws.Range("B17").Select()
